# Correção do slide sobre manipulação de estados
# Slide 15 ("No próximo módulo vamos aprender sobre armazenamento
# interno com arquivos.") needs to become bold.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(1)          # "TextBox 13"
$tr = $shp.TextFrame.TextRange

$tr.Font.Bold = 1
